$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the tail of the header row right by two columns (AE:AJ -> AG:AL is NOT what happens;
# instead we directly retype each header cell in its new position), inserting two new
# headers ("Product", "API") right after "Quantity4", and collapsing the four
# "Exit Price1".."Exit Price4" columns into a single renamed "CurrentSop" column.

$ws.Range("AC1").Value() = "Product"
$ws.Range("AD1").Value() = "API"
$ws.Range("AE1").Value() = "count"
$ws.Range("AF1").Value() = "isop"
$ws.Range("AG1").Value() = "Completed Date"
$ws.Range("AH1").Value() = "CurrentSop"

# Drop the now-unused trailing columns (old Exit Price1 / Exit Price2 / Exit Price3 / Exit Price4)
$ws.Range("AI1:AJ1").ClearContents()

# The two newly-introduced header cells pick up the same fill formatting as their
# neighbour (Quantity4 / AB1)
$ws.Range("AB1").Copy()
$ws.Range("AC1:AD1").PasteSpecial(-4122)
$ws.Range("AC1").Value() = "Product"
$ws.Range("AD1").Value() = "API"

# Widen the two columns that now hold "isop" / "Completed Date" so the longer labels fit
$ws.Columns("AF:AF").ColumnWidth = 11.65
$ws.Columns("AG:AG").ColumnWidth = 13.8

$ws.Range("AC6").Select()
